# Update the odds data on the active sheet (rows 3, 5, 6) to reflect
# refreshed FlashScore odds values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 1.6
$ws.Range("J3").Value = 5.5
$ws.Range("K3").Value = 2.12
$ws.Range("L3").Value = 2.18
$ws.Range("N3").Value = 7.2
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.2
$ws.Range("Q3").Value = 1.91
$ws.Range("R3").Value = 1.82
$ws.Range("S3").Value = 1.42
$ws.Range("T3").Value = 2.65
$ws.Range("U3").Value = 1.93
$ws.Range("V3").Value = 1.78
$ws.Range("W3").Value = 14
$ws.Range("Y3").Value = 17.5
$ws.Range("Z3").Value = 110
$ws.Range("AA3").Value = 60
$ws.Range("AB3").Value = 60
$ws.Range("AC3").Value = 7.2
$ws.Range("AD3").Value = 7.2
$ws.Range("AE3").Value = 17.5
$ws.Range("AF3").Value = 90
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 6.3
$ws.Range("AI3").Value = 7
$ws.Range("AM3").Value = 28
$ws.Range("AO3").Value = 32
$ws.Range("AP3").Value = 37
$ws.Range("AR3").Value = 250
$ws.Range("AT3").Value = 2.65
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 80
$ws.Range("AW3").Value = 3.35
$ws.Range("AX3").Value = 7.9
$ws.Range("AY3").Value = 19
$ws.Range("AZ3").Value = 26
$ws.Range("BA3").Value = 65
$ws.Range("BB3").Value = 300

# --- Row 5 ---
$ws.Range("K5").Value = 2.38
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.25
$ws.Range("W5").Value = 7.5
$ws.Range("AT5").Value = 3.25
$ws.Range("AW5").Value = 7.5

# --- Row 6 ---
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 8.5
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.85
$ws.Range("R6").Value = 1.95
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("AH6").Value = 6.5
$ws.Range("AT6").Value = 3
